$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# Simple one-for-one text replacements (single run, no structural change)
Replace-Exact "What is a focus group and what will this one look like?" "Yintoni iqela ekugxilwe kulo kwaye eli liza kujongeka njani?"

Replace-Exact "Why have I been invited to the interview?" "Kutheni ndimenyiwe nje kudliwano-ndlebe?"

Replace-Exact "Do I get anything for being interviewed? " "Ingaba ikhona into endiyifumanayo ngokwenziwa oludliwano-ndlebe? "

Replace-Exact "As a thank you for taking part in the discussion, we'll give you a R120 Shoprite voucher afterwards. " "Njengombulelo ngokuthatha inxaxheba kwingxoxo, siza kukunika ivawutsha ye-R120 yakwaShoprite emva koko."

Replace-Exact "Who pays for the study?" "Ngubani obhatalela oluphononongo?"

Replace-Exact "This study is part of the Global Parenting Initiative, funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. " "Olu phononongo luyinxalenye ye Global Parenting Initiative, luxhaswe ngokwezimali ngu LEGO Foundation, Oak Fundation, i-World Childhood Foundation, i-Human Safety Net kunye ne UK Research kunye ne Innovaion Global Challenges Research Fund. "

Replace-Exact "Data protection" "Ukhuseleko lwedatha"

Replace-Exact "I know who can see my information after the focus group, how it will be kept safe, and what happens to it after the study." "Ndiyayazi ukuba ngubani onokubona ulwazi lwam emva kweqela ekugxilwe kulo, ukuba luya kugcinwa njani lukhuselekile, kwaye kwenzeka ntoni kulo emva kophononongo."

# Cases that also drop a trailing run (a lone space/nbsp run) -- the whole
# paragraph (minus its end-of-paragraph mark) is replaced with the new text
# as a single run, inheriting the formatting of the paragraph's first run.
function Replace-ParagraphByText($oldStart, $newText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.StartsWith($oldStart)) {
            $r = $p.Range
            $r2 = $d.Range($r.Start, $r.End - 1)
            $r2.Text = $newText
            return
        }
    }
    Write-Output "PARAGRAPH NOT FOUND: $oldStart"
}

Replace-ParagraphByText "What happens to my information if I agree to be interviewed?" "Kwenzeka ntoni ngeenkcukacha zam ukuba ndiyavuma ukuba noludliwano-ndlebe?"

Replace-ParagraphByText "Who are some of the study team members?" "Ngobani amanye amalungu eqela lophononongo?"

# Remaining simple replacement
Replace-Exact "The principal investigators of this study are Prof Cathy Ward and Cindee Bruyns and the Co-investigator is Carly Katzef all from the University of Cape Town." "Abaphononongi abaziintloko kolu phononongo nguNjinga Cathy Ward no Cindee Bruyns ze Co-investigator ngu Carly Katzef bonke basuka kwiDyunivesithi yaseKapa."
